# Add "Disease Mapping 2" (module 5) vocabulary terms to the glossary, and
# touch up the formatting on the first new definition cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$terms = @(
    @("Distance", "A fundamental dimension in geography referring to the strength of connectedness or proximity in eculidean space, social space, or network space. Distance if fundamental because we assume that a) entities that are closer are, on average, more alike than entities that are far apart; and b) increasing distance represents increasing friction or imedance to social and health-relevant interaction"),
    @("Aspatial vs. Spatial", "This distinction refers to whether or not spatial proximity or contiguity is explicitly incorporated into an analysis (spatial) versus whether spatial units are treated as independent of one another (aspatial)"),
    @("Spatial neighbors", "The set of spatial entities that are determined to be 'near' rather than 'far' (in binary terms) or relatively 'closer' or 'further' (in continuous terms). The definition of 'neighbors' is part of specifying spatial relatedness."),
    @("Spatial weights matrix", "Typically a square matrix (n rows x n columns where n=geographic units) indexing all units on rows and columns. The values in the matrix indicate the spatial connectedness between all pairs of units."),
    @("Neighbor symmetry", "An attribute of spatial relationships in which it is assumed that if spatial unit A is a neighbor with B, then spatial unit B is also a neighbor with A. Some neighbor definitions (e.g. k-nearest neighbors) do not require symmetry."),
    @("Toblers' First Law of Geography", "All things are related, but near things are more related on average than distant things"),
    @("Delauney triangulation", "Geometric strategy for creating a mesh of contiguous, nonoverlapping triangles from a dataset of points. If points are the centroids of polygons, the triangle edges become graph-based definitions of spatial neighbors")
)

$row = 31
foreach ($pair in $terms) {
    $ws.Cells.Item($row, 1).Value = 5
    $ws.Cells.Item($row, 2).Value = $pair[0]
    $ws.Cells.Item($row, 3).Value = $pair[1]
    $row++
}

# Small formatting fix touching the first new definition cell (creates the
# new cellXfs entry with the alignment flag applied).
$ws.Range("C31").WrapText = $true

# Scroll/selection state left over from editing near the bottom of the sheet.
$ws.Activate()
$ws.Range("A28").Select()
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 1
